$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9; this shifts the existing rows 9-12 down to 10-13,
# carrying their values and formatting (including the date-format style on column D).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44489
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100101
$ws.Cells.Item(9, 8).Value = "Berries"
$ws.Cells.Item(9, 9).Value = 100101001
$ws.Cells.Item(9, 10).Value = "Arándano (blue)"
$ws.Cells.Item(9, 11).Value = "Sin especificar"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 400
$ws.Cells.Item(9, 14).Value = 11500
$ws.Cells.Item(9, 15).Value = 12000
$ws.Cells.Item(9, 16).Value = 11750
$ws.Cells.Item(9, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 19).Value = 5875
$ws.Cells.Item(9, 20).Value = 2
